# Update "Generate Report for Handback" timestamps in the handback-status workbook.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G3 - Latest HO Xliff Generate Date for the 3a6469dd... file.
# (This text is shared with de-de!H3, so both update together.)
$wsOverview.Range("G3").Value = "2016-08-28 22:45:30"

# zh-cn!H3 - Correspond Handoff Datetime for the 3a6469dd... file.
$wsZhCn.Range("H3").Value = "2016-08-28 22:45:26"

# zh-cn!K3 - Correspond Handback DateTime for the 3a6469dd... file.
$wsZhCn.Range("K3").Value = "2016-08-28 22:45:43"

# de-de!H3 - Correspond Handoff Datetime for the 3a6469dd... file.
# (Shares the same underlying text as Overview!G3.)
$wsDeDe.Range("H3").Value = "2016-08-28 22:45:30"

# de-de!K3 - Correspond Handback DateTime for the 3a6469dd... file.
$wsDeDe.Range("K3").Value = "2016-08-28 22:45:51"
